# Apply the MASTERSHEET.xlsx changes described in the commit
# "mergin with Nitish's changes":
#  - C2 (EXECUTE flag): Yes -> No
#  - B20 (Test Case Description): drop the "(Flow7)" suffix
#  - C20 (EXECUTE flag): No -> Yes
#  - B21 (Test Case Description): drop the "(Flow3)" suffix
#  - C21 (EXECUTE flag): No -> Yes
#  - B22 (Test Case Description): drop the "(Flow5)" suffix
#  - move the active selection to C22

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "No"

$ws.Range("B20").Value = "Joint Account for One New and One Existing Customer"
$ws.Range("C20").Value = "Yes"

$ws.Range("B21").Value = "Individual Account for Existing Customer having Joint Account Only"
$ws.Range("C21").Value = "Yes"

$ws.Range("B22").Value = "Joint Account for both Existing Customers"

$ws.Range("C22").Select()
